# Updates the calibration figures on the AR / SETAR / GARCH / TARCH / AR_TARCH
# sheets of the gold riskDriverType PnL factor calibration workbook
# ("Finalized definition of virtual environment").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# AR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.0008508143752312533
$ws.Range("B3").Value = 0.7827615911638615
$ws.Range("B4").Value = 0.0128064591397295
$ws.Range("B5").Value = "[0.9999999999999999, 0.26106732725587745, 0.21280658769846864, 0.18530453660527915, 0.19358804061053356, 0.38822353941157006, 0.11753214282261931, 0.09628826260800163, 0.04107700761261552, 0.06740511094199222, 2.696898192855705e-05, 0.0007064898455680466, 0.008349365799359344, 0.010093739824084317, 0.03705100645908723, -0.008397211666626674, -0.0017107622328299886, -0.01431526915190963, 0.03313045865462108, 0.017641536576249135]"

# ---------------------------------------------------------------------------
# SETAR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B3").Value = 0.7985294117647059
$ws.Range("B4").Value = -0.0323797923443367
$ws.Range("B5").Value = 0.7401583960011209
$ws.Range("B6").Value = 0.01002174191600099
$ws.Range("B7").Value = 0.02911954184775502
$ws.Range("B8").Value = 0.7257135347815349
$ws.Range("B9").Value = 0.008551516766166438
$ws.Range("B10").Value = "[1.0, 0.2327024025292807, 0.21122730975867723, 0.17725693372802087, 0.14523747743603269, 0.273420657615024, 0.11294257557126536, 0.12748543767261059, 0.08118207967815194, 0.06734918909274266, 0.03947556012582784, 0.051191439495166906, 0.04865601534863883, 0.053451249240821404, 0.07653074041577203, 0.040378220947284424, 0.05400854688614327, 0.018628761892072156, 0.03851933777439422, 0.033597445210607736]"

# ---------------------------------------------------------------------------
# GARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.0001418582543074643
$ws.Range("B3").Value = 0.0005940176682119875
$ws.Range("B4").Value = 0.08045190067059058
$ws.Range("B5").Value = 0.8780386739425166
$ws.Range("B6").Value = "[1.0, 0.02573486851203955, 0.0015710839042032485, -0.022769811650114418, -0.030112163396515196, 0.35706895948981654, 0.014763016189585013, 0.016767749287065883, -0.031040846760766332, -0.0314173078087463, -0.07786642283245827, -0.07209846803747898, -0.0585581529696338, -0.01395421147950501, 0.0015540352602608969, -0.04728804982036745, -0.03955523243373323, -0.07912840478039397, -0.020209928991535484, -0.029570264828236985]"

# ---------------------------------------------------------------------------
# TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.001702897179315271
$ws.Range("B3").Value = 0.0005562539495772171
$ws.Range("B4").Value = 0.058996798770054
$ws.Range("B5").Value = 0.8840047700633431
$ws.Range("B6").Value = "[0.9999999999999999, 0.026422080136563483, -0.0005556620091074283, -0.024115094362629905, -0.03000326156467007, 0.35753942478363326, 0.015716435573832236, 0.01605491278111625, -0.03156720915181818, -0.030713317251441732, -0.0765002370066913, -0.07318333794870993, -0.05915773884003885, -0.014618554612888827, 0.0005497908719769208, -0.04752797410745079, -0.039693915410176764, -0.07919951753320485, -0.02006315725933448, -0.029908405187217903]"
$ws.Range("B7").Value = 0.03664180900297866

# ---------------------------------------------------------------------------
# AR_TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.0006554969042932606
$ws.Range("B3").Value = 0.0007633141250623466
$ws.Range("B4").Value = 0.07375687319921803
$ws.Range("B5").Value = 0.8486602057783299
$ws.Range("B6").Value = "[1.0, 0.09724507159753526, 0.07181833348276442, 0.052412858121071576, 0.04803503060565373, 0.256129136937929, -0.0013701382589351536, -0.0005782933448187097, -0.04702468791079992, -0.02351109474500065, -0.07807620111662228, -0.0640318214012896, -0.05503427631844634, -0.034354139225936164, -0.009353210029520788, -0.05610617799292541, -0.046238704806255644, -0.05616256541811161, -0.012864277204530346, -0.022790598152768838]"
$ws.Range("B7").Value = 0.02843397606218394
$ws.Range("B9").Value = 0.7524462657708358
